$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the simulation log values (run_time, max_er, iter 0..19) for the
# gr75_01 dataset rows after switching the log write mode.
    $ws.Range("C2").Value = 0.5269255638122559
    $ws.Range("E2").Value = 603.0102269504332
    $ws.Range("F2").Value = 0.02113855288568531
    $ws.Range("G2").Value = 0.01816047106509355
    $ws.Range("H2").Value = 0.01668346413106496
    $ws.Range("I2").Value = 0.01576607470620122
    $ws.Range("J2").Value = 0.01472376662218459
    $ws.Range("K2").Value = 0.01384707527371461
    $ws.Range("L2").Value = 0.01338806196857388
    $ws.Range("M2").Value = 0.01289018853806115
    $ws.Range("N2").Value = 0.01289018853806115
    $ws.Range("O2").Value = 0.0127608088633898
    $ws.Range("P2").Value = 0.0127608088633898
    $ws.Range("Q2").Value = 0.01253698328309324
    $ws.Range("R2").Value = 0.01238062877748528
    $ws.Range("S2").Value = 0.01226850010260345
    $ws.Range("T2").Value = 0.01208062205815521
    $ws.Range("U2").Value = 0.01195082063860927
    $ws.Range("V2").Value = 0.01193258432794835
    $ws.Range("W2").Value = 0.01183970457068713
    $ws.Range("X2").Value = 0.01179986567307338
    $ws.Range("Y2").Value = 0.01175458532067121
    $ws.Range("C3").Value = 0.5312392711639404
    $ws.Range("E3").Value = 588.061798362929
    $ws.Range("F3").Value = 0.02090514152086939
    $ws.Range("G3").Value = 0.01800525247790521
    $ws.Range("H3").Value = 0.01547576836219271
    $ws.Range("I3").Value = 0.0148512500673906
    $ws.Range("J3").Value = 0.01385570876662849
    $ws.Range("K3").Value = 0.01340208452406244
    $ws.Range("L3").Value = 0.01306155858872731
    $ws.Range("M3").Value = 0.01264463966712405
    $ws.Range("N3").Value = 0.01225395508262999
    $ws.Range("O3").Value = 0.01219996966205949
    $ws.Range("P3").Value = 0.01197793934936026
    $ws.Range("Q3").Value = 0.01173498110618846
    $ws.Range("R3").Value = 0.01173498110618846
    $ws.Range("S3").Value = 0.01168121062535062
    $ws.Range("T3").Value = 0.01161870411830767
    $ws.Range("U3").Value = 0.01157360523436896
    $ws.Range("V3").Value = 0.01154468307238202
    $ws.Range("W3").Value = 0.01151171009523183
    $ws.Range("X3").Value = 0.01148016090733078
    $ws.Range("Y3").Value = 0.01146319295054442
    $ws.Range("C4").Value = 0.6441624164581299
    $ws.Range("E4").Value = 598.2268122965252
    $ws.Range("F4").Value = 0.02145118859480475
    $ws.Range("G4").Value = 0.01853430679924043
    $ws.Range("H4").Value = 0.01660950980965648
    $ws.Range("I4").Value = 0.01565832752383478
    $ws.Range("J4").Value = 0.01445408703445077
    $ws.Range("K4").Value = 0.01373455404343684
    $ws.Range("L4").Value = 0.0131310700741944
    $ws.Range("M4").Value = 0.0131310700741944
    $ws.Range("N4").Value = 0.01301432760454357
    $ws.Range("O4").Value = 0.01280479219649386
    $ws.Range("P4").Value = 0.01253796721413884
    $ws.Range("Q4").Value = 0.01231572641287818
    $ws.Range("R4").Value = 0.01219278675906894
    $ws.Range("S4").Value = 0.01194377400292806
    $ws.Range("T4").Value = 0.01189531210093814
    $ws.Range("U4").Value = 0.01179396738180705
    $ws.Range("V4").Value = 0.01176528877973847
    $ws.Range("W4").Value = 0.01171967478834992
    $ws.Range("X4").Value = 0.01168293416572071
    $ws.Range("Y4").Value = 0.01166134137030263
    $ws.Range("C5").Value = 0.5942234992980957
    $ws.Range("E5").Value = 602.3546427493802
    $ws.Range("F5").Value = 0.02123306595234294
    $ws.Range("G5").Value = 0.01824632604175058
    $ws.Range("H5").Value = 0.01627013399764592
    $ws.Range("I5").Value = 0.0157782434344335
    $ws.Range("J5").Value = 0.01506721532839755
    $ws.Range("K5").Value = 0.01441304778509578
    $ws.Range("L5").Value = 0.01359608201450409
    $ws.Range("M5").Value = 0.01297699744629068
    $ws.Range("N5").Value = 0.0129018261793614
    $ws.Range("O5").Value = 0.01280821421680652
    $ws.Range("P5").Value = 0.01266951072788022
    $ws.Range("Q5").Value = 0.01226315548250724
    $ws.Range("R5").Value = 0.01224575403033537
    $ws.Range("S5").Value = 0.01207131687933331
    $ws.Range("T5").Value = 0.01200971557112397
    $ws.Range("U5").Value = 0.01195719961133055
    $ws.Range("V5").Value = 0.01189423606163728
    $ws.Range("W5").Value = 0.01177787364022701
    $ws.Range("X5").Value = 0.01177434862503124
    $ws.Range("Y5").Value = 0.01174180590154737
    $ws.Range("C6").Value = 0.5468759536743164
    $ws.Range("E6").Value = 597.1654948600299
    $ws.Range("F6").Value = 0.02109460741991599
    $ws.Range("G6").Value = 0.01747579612791106
    $ws.Range("H6").Value = 0.01620861065550056
    $ws.Range("I6").Value = 0.01561902092230849
    $ws.Range("J6").Value = 0.01437192966169548
    $ws.Range("K6").Value = 0.01407170266500568
    $ws.Range("L6").Value = 0.01354735257671879
    $ws.Range("M6").Value = 0.01335378399284603
    $ws.Range("N6").Value = 0.01275197146844777
    $ws.Range("O6").Value = 0.01263295534184654
    $ws.Range("P6").Value = 0.01243948001027791
    $ws.Range("Q6").Value = 0.01231304638661348
    $ws.Range("R6").Value = 0.01215421027852838
    $ws.Range("S6").Value = 0.01205465435362124
    $ws.Range("T6").Value = 0.01198686151109865
    $ws.Range("U6").Value = 0.01193088975993096
    $ws.Range("V6").Value = 0.01182946698487079
    $ws.Range("W6").Value = 0.01170297871550362
    $ws.Range("X6").Value = 0.01167481383402817
    $ws.Range("Y6").Value = 0.01164065292124814
    $ws.Range("C7").Value = 0.5746200084686279
    $ws.Range("E7").Value = 608.2545791985103
    $ws.Range("F7").Value = 0.02133924811251438
    $ws.Range("G7").Value = 0.01845199019616885
    $ws.Range("H7").Value = 0.01594833130070882
    $ws.Range("I7").Value = 0.01501344565282738
    $ws.Range("J7").Value = 0.01444206089325096
    $ws.Range("K7").Value = 0.01404820424223
    $ws.Range("L7").Value = 0.01361442005576446
    $ws.Range("M7").Value = 0.01313415184064424
    $ws.Range("N7").Value = 0.01275323565336737
    $ws.Range("O7").Value = 0.01266388032346212
    $ws.Range("P7").Value = 0.01261142366286919
    $ws.Range("Q7").Value = 0.0124414011602778
    $ws.Range("R7").Value = 0.01236973190614907
    $ws.Range("S7").Value = 0.01224107387420039
    $ws.Range("T7").Value = 0.01209611598739189
    $ws.Range("U7").Value = 0.01204784424768465
    $ws.Range("V7").Value = 0.01201214979314516
    $ws.Range("W7").Value = 0.01192248016603698
    $ws.Range("X7").Value = 0.01192248016603698
    $ws.Range("Y7").Value = 0.01185681440932768
    $ws.Range("C8").Value = 0.5312502384185791
    $ws.Range("E8").Value = 602.9935630297841
    $ws.Range("F8").Value = 0.02120117246638004
    $ws.Range("G8").Value = 0.01808265187818603
    $ws.Range("H8").Value = 0.01636532109824761
    $ws.Range("I8").Value = 0.01539899677418049
    $ws.Range("J8").Value = 0.01482301170303223
    $ws.Range("K8").Value = 0.01426466477160838
    $ws.Range("L8").Value = 0.01334343795620644
    $ws.Range("M8").Value = 0.01287781567810648
    $ws.Range("N8").Value = 0.01287781567810648
    $ws.Range("O8").Value = 0.01263650010025517
    $ws.Range("P8").Value = 0.01243008358516693
    $ws.Range("Q8").Value = 0.0123997940297162
    $ws.Range("R8").Value = 0.01232530465439132
    $ws.Range("S8").Value = 0.01205838002392063
    $ws.Range("T8").Value = 0.01203863691268674
    $ws.Range("U8").Value = 0.01195078747779424
    $ws.Range("V8").Value = 0.01188539998353398
    $ws.Range("W8").Value = 0.01183143080469344
    $ws.Range("X8").Value = 0.0117905826385242
    $ws.Range("Y8").Value = 0.01175426048791002
    $ws.Range("C9").Value = 0.5312364101409912
    $ws.Range("E9").Value = 603.7738560027719
    $ws.Range("F9").Value = 0.02123788799319374
    $ws.Range("G9").Value = 0.01779393992738511
    $ws.Range("H9").Value = 0.01655167779717226
    $ws.Range("I9").Value = 0.01561977342560152
    $ws.Range("J9").Value = 0.01449936617384391
    $ws.Range("K9").Value = 0.01403846960159194
    $ws.Range("L9").Value = 0.01341656826271438
    $ws.Range("M9").Value = 0.01332291727329071
    $ws.Range("N9").Value = 0.01310259677129597
    $ws.Range("O9").Value = 0.01283556905408317
    $ws.Range("P9").Value = 0.01263838196570119
    $ws.Range("Q9").Value = 0.01236933473522074
    $ws.Range("R9").Value = 0.01233076175442323
    $ws.Range("S9").Value = 0.0121166156369335
    $ws.Range("T9").Value = 0.01204945266813324
    $ws.Range("U9").Value = 0.01191665264565152
    $ws.Range("V9").Value = 0.01189591256246675
    $ws.Range("W9").Value = 0.01184213544632043
    $ws.Range("X9").Value = 0.01180033298420681
    $ws.Range("Y9").Value = 0.01176947087724701
    $ws.Range("C10").Value = 0.5468900203704834
    $ws.Range("E10").Value = 625.3663212767151
    $ws.Range("F10").Value = 0.02134107648440884
    $ws.Range("G10").Value = 0.01837935973675144
    $ws.Range("H10").Value = 0.01625816187309214
    $ws.Range("I10").Value = 0.01561396340321121
    $ws.Range("J10").Value = 0.01444468405749754
    $ws.Range("K10").Value = 0.01431280879799099
    $ws.Range("L10").Value = 0.01391984862407272
    $ws.Range("M10").Value = 0.01364358958402466
    $ws.Range("N10").Value = 0.01364358958402466
    $ws.Range("O10").Value = 0.01336269938969645
    $ws.Range("P10").Value = 0.01296934559670344
    $ws.Range("Q10").Value = 0.01289180516903879
    $ws.Range("R10").Value = 0.0127089101247955
    $ws.Range("S10").Value = 0.01260525551127072
    $ws.Range("T10").Value = 0.01249983126409438
    $ws.Range("U10").Value = 0.01247193758261591
    $ws.Range("V10").Value = 0.012361520118612
    $ws.Range("W10").Value = 0.01226730431013635
    $ws.Range("X10").Value = 0.01222260344062877
    $ws.Range("Y10").Value = 0.01219037663307436
    $ws.Range("C11").Value = 0.5312352180480957
    $ws.Range("E11").Value = 599.4369188301825
    $ws.Range("F11").Value = 0.02130631467251787
    $ws.Range("G11").Value = 0.01798507629329489
    $ws.Range("H11").Value = 0.01614852946204287
    $ws.Range("I11").Value = 0.01541448010276658
    $ws.Range("J11").Value = 0.01429049516936974
    $ws.Range("K11").Value = 0.01391044455987069
    $ws.Range("L11").Value = 0.01331305916599156
    $ws.Range("M11").Value = 0.01312464654957834
    $ws.Range("N11").Value = 0.01259924854687428
    $ws.Range("O11").Value = 0.01259802147006913
    $ws.Range("P11").Value = 0.01236496640629553
    $ws.Range("Q11").Value = 0.01226947145725192
    $ws.Range("R11").Value = 0.01220872863042655
    $ws.Range("S11").Value = 0.01200739242156962
    $ws.Range("T11").Value = 0.01189417077853528
    $ws.Range("U11").Value = 0.01184675837383901
    $ws.Range("V11").Value = 0.01178507711936227
    $ws.Range("W11").Value = 0.01172538285504021
    $ws.Range("X11").Value = 0.01171723878156945
    $ws.Range("Y11").Value = 0.01168493019162149
